# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job Leve
# tables. Values below were recomputed from refreshed source data; a few
# rows had their NQ/HQ split collapse back to a single (NQ-only) price,
# which drops the HQ-profit figure (column N) for those rows entirely.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H4").Value = 314.44446
$ws.Range("I4").Value = 314.44446
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 314.44446
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -200.44446
$ws.Range("N4").ClearContents()

$ws.Range("H39").Value = 557.3
$ws.Range("I39").Value = 71.63636
$ws.Range("J39").Value = 1150.8889
$ws.Range("K39").Value = 214.90908
$ws.Range("L39").Value = 3452.6667
$ws.Range("M39").Value = 81.09092000000001
$ws.Range("N39").Value = -4044.6667

$ws.Range("H106").Value = 900
$ws.Range("I106").Value = 920
$ws.Range("J106").Value = 800
$ws.Range("K106").Value = 920
$ws.Range("L106").Value = 800
$ws.Range("M106").Value = -289
$ws.Range("N106").Value = -2062

$ws.Range("H137").Value = 3229027.8
$ws.Range("I137").Value = 9096400
$ws.Range("J137").Value = 1973.1
$ws.Range("K137").Value = 27289200
$ws.Range("L137").Value = 5919.299999999999
$ws.Range("M137").Value = -27286650
$ws.Range("N137").Value = -11019.3

$ws.Range("H141").Value = 589799.9399999999
$ws.Range("I141").Value = 1634.6666
$ws.Range("J141").Value = 2060213.1
$ws.Range("K141").Value = 4903.9998
$ws.Range("L141").Value = 6180639.300000001
$ws.Range("M141").Value = 276.0002000000004
$ws.Range("N141").Value = -6190999.300000001

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H54").Value = 2860.5
$ws.Range("I54").Value = 2860.5
$ws.Range("K54").Value = 2860.5
$ws.Range("M54").Value = -2376.5

$ws.Range("H99").Value = 1700.7778
$ws.Range("I99").Value = 1260.5
$ws.Range("J99").Value = 2958.7144
$ws.Range("K99").Value = 1260.5
$ws.Range("L99").Value = 2958.7144
$ws.Range("M99").Value = 237.5
$ws.Range("N99").Value = -5954.7144

$ws.Range("H107").Value = 2030.5834
$ws.Range("I107").Value = 1247.2858
$ws.Range("J107").Value = 3127.2
$ws.Range("K107").Value = 1247.2858
$ws.Range("L107").Value = 3127.2
$ws.Range("M107").Value = 672.7141999999999
$ws.Range("N107").Value = -6967.2

$ws.Range("H134").Value = 2850.6511
$ws.Range("I134").Value = 2566.5144
$ws.Range("J134").Value = 4093.75
$ws.Range("K134").Value = 7699.5432
$ws.Range("L134").Value = 12281.25
$ws.Range("M134").Value = -5164.5432
$ws.Range("N134").Value = -17351.25

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H31").Value = 1236531.4
$ws.Range("I31").Value = 2041970.1
$ws.Range("J31").Value = 3203.4688
$ws.Range("K31").Value = 2041970.1
$ws.Range("L31").Value = 3203.4688
$ws.Range("M31").Value = -2041675.1
$ws.Range("N31").Value = -3793.4688

$ws.Range("H34").Value = 1236531.4
$ws.Range("I34").Value = 2041970.1
$ws.Range("J34").Value = 3203.4688
$ws.Range("K34").Value = 2041970.1
$ws.Range("L34").Value = 3203.4688
$ws.Range("M34").Value = -2041768.1
$ws.Range("N34").Value = -3607.4688

$ws.Range("H48").Value = 54788.25
$ws.Range("J48").Value = 54788.25
$ws.Range("L48").Value = 54788.25
$ws.Range("N48").Value = -55740.25

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 3583.5
$ws.Range("I132").Value = 2441.25
$ws.Range("J132").Value = 4725.75
$ws.Range("K132").Value = 7323.75
$ws.Range("L132").Value = 14177.25
$ws.Range("M132").Value = -4793.75
$ws.Range("N132").Value = -19237.25

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 1744.8223
$ws.Range("I134").Value = 1407.3721
$ws.Range("J134").Value = 9000
$ws.Range("K134").Value = 4222.1163
$ws.Range("L134").Value = 27000
$ws.Range("M134").Value = -1687.1163
$ws.Range("N134").Value = -32070

$ws.Range("H137").Value = 15000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H68").Value = 1821.9642
$ws.Range("I68").Value = 620.46155
$ws.Range("J68").Value = 2863.2666
$ws.Range("K68").Value = 1861.38465
$ws.Range("L68").Value = 8589.799800000001
$ws.Range("M68").Value = -1050.38465
$ws.Range("N68").Value = -10211.7998

$ws.Range("H71").Value = 1821.9642
$ws.Range("I71").Value = 620.46155
$ws.Range("J71").Value = 2863.2666
$ws.Range("K71").Value = 5584.15395
$ws.Range("L71").Value = 25769.3994
$ws.Range("M71").Value = -1528.15395
$ws.Range("N71").Value = -33881.39939999999

$ws.Range("H107").Value = 834.8570999999999
$ws.Range("I107").Value = 546.7879
$ws.Range("J107").Value = 1151.7333
$ws.Range("K107").Value = 1640.3637
$ws.Range("L107").Value = 3455.199900000001
$ws.Range("M107").Value = 279.6362999999999
$ws.Range("N107").Value = -7295.199900000001

$ws.Range("H122").Value = 897.5294
$ws.Range("J122").Value = 1917.3334
$ws.Range("L122").Value = 17256.0006
$ws.Range("N122").Value = -22156.0006

$ws.Range("H131").Value = 1410.9814
$ws.Range("J131").Value = 1068.8536
$ws.Range("L131").Value = 3206.5608
$ws.Range("N131").Value = -13286.5608

$ws.Range("H132").Value = 2872
$ws.Range("I132").Value = 2319.4
$ws.Range("J132").Value = 3332.5
$ws.Range("K132").Value = 20874.6
$ws.Range("L132").Value = 29992.5
$ws.Range("M132").Value = -18344.6
$ws.Range("N132").Value = -35052.5

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H132").Value = 3304.6445
$ws.Range("I132").Value = 2413.5806
$ws.Range("K132").Value = 7240.7418
$ws.Range("M132").Value = -4710.7418

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 200002960
$ws.Range("I22").Value = 500000500
$ws.Range("J22").Value = 4600
$ws.Range("K22").Value = 500000500
$ws.Range("L22").Value = 4600
$ws.Range("M22").Value = -500000205
$ws.Range("N22").Value = -5190

$ws.Range("H27").Value = 200002960
$ws.Range("I27").Value = 500000500
$ws.Range("J27").Value = 4600
$ws.Range("K27").Value = 500000500
$ws.Range("L27").Value = 4600
$ws.Range("M27").Value = -500000393
$ws.Range("N27").Value = -4814

$ws.Range("H68").Value = 1428.5714
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 1000
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -251
$ws.Range("N68").Value = -11498

$ws.Range("H71").Value = 1428.5714
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 5000
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -1256
$ws.Range("N71").Value = -57488

$ws.Range("H82").Value = 2674.52
$ws.Range("I82").Value = 1993.0769
$ws.Range("J82").Value = 3412.75
$ws.Range("K82").Value = 1993.0769
$ws.Range("L82").Value = 3412.75
$ws.Range("M82").Value = -1632.0769
$ws.Range("N82").Value = -4134.75

$ws.Range("H85").Value = 2674.52
$ws.Range("I85").Value = 1993.0769
$ws.Range("J85").Value = 3412.75
$ws.Range("K85").Value = 1993.0769
$ws.Range("L85").Value = 3412.75
$ws.Range("M85").Value = -745.0769
$ws.Range("N85").Value = -5908.75

$ws.Range("H132").Value = 3200.25
$ws.Range("I132").Value = 2245.7273
$ws.Range("J132").Value = 4007.923
$ws.Range("K132").Value = 6737.1819
$ws.Range("L132").Value = 12023.769
$ws.Range("M132").Value = -4207.1819
$ws.Range("N132").Value = -17083.769

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H132").Value = 200444.08
$ws.Range("I132").Value = 295889.56
$ws.Range("J132").Value = 9553.117
$ws.Range("K132").Value = 887668.6799999999
$ws.Range("L132").Value = 28659.351
$ws.Range("M132").Value = -885138.6799999999
$ws.Range("N132").Value = -33719.351

$ws.Range("H136").Value = 1417.1305
$ws.Range("I136").Value = 673.4211
$ws.Range("J136").Value = 4949.75
$ws.Range("K136").Value = 2020.2633
$ws.Range("L136").Value = 14849.25
$ws.Range("M136").Value = 529.7366999999999
$ws.Range("N136").Value = -19949.25

